# Applies the "Updated cryptos list" refresh: new Price (column D) and
# Volume(1h) (column E) values for the affected rows of the crypto table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.512.12"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "1.621.36"
$ws.Range("E3").Value = "  -1.42%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'211.44"  # stored as text, same as original
$ws.Range("E5").Value = "  -0.78%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.41%  "
$ws.Range("D9").Value = "'0.263"  # stored as text, same as original
$ws.Range("E9").Value = "  +1.62%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("D11").Value = "'0.0881"  # stored as text, same as original
$ws.Range("E11").Value = "  -1.50%  "
$ws.Range("D12").Value = "1.850.76"
$ws.Range("E12").Value = "  -1.45%  "
$ws.Range("D13").Value = "1.636.76"
$ws.Range("E13").Value = "  -0.35%  "
$ws.Range("E14").Value = "  -0.19%  "
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").Value = "'65.20"  # stored as text, same as original
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "27.489.83"
$ws.Range("D18").Value = "'229.52"  # stored as text, same as original
$ws.Range("E18").Value = "  -0.86%  "
$ws.Range("E19").Value = "  -0.93%  "
$ws.Range("E20").Value = "  -2.09%  "
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "'2.11"  # stored as text, same as original
$ws.Range("E24").Value = "  +7.39%  "
$ws.Range("D25").Value = "'149.35"  # stored as text, same as original
$ws.Range("E25").Value = "  -0.46%  "
$ws.Range("E26").Value = "  -1.15%  "
$ws.Range("E27").Value = "  -0.88%  "
$ws.Range("E28").Value = "  -0.03%  "
$ws.Range("D29").Value = "'15.53"  # stored as text, same as original
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'1.17"  # stored as text, same as original
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("E32").Value = "  -1.26%  "
$ws.Range("D33").Value = "1.463.99"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("E34").Value = "  -2.85%  "
$ws.Range("E35").Value = "  -2.38%  "
$ws.Range("E36").Value = "  -0.40%  "
$ws.Range("D37").Value = "'0.941"  # stored as text, same as original
$ws.Range("E37").Value = "  +6.30%  "
$ws.Range("E38").Value = "  -0.37%  "
$ws.Range("D39").Value = "'0.871"  # stored as text, same as original
$ws.Range("E39").Value = "  -1.54%  "
$ws.Range("D40").Value = "'0.551"  # stored as text, same as original
$ws.Range("E40").Value = "  -3.32%  "
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").Value = "'67.33"  # stored as text, same as original
$ws.Range("E43").Value = "  -5.36%  "
$ws.Range("E44").Value = "  -0.74%  "
$ws.Range("E45").Value = "  -2.09%  "
$ws.Range("E46").Value = "  -6.17%  "
$ws.Range("D47").Value = "'1.75"  # stored as text, same as original
$ws.Range("E47").Value = "  +1.48%  "
$ws.Range("D48").Value = "1.761.05"
$ws.Range("E48").Value = "  -1.51%  "
$ws.Range("E49").Value = "  +1.54%  "
$ws.Range("E50").Value = "  -1.00%  "
$ws.Range("E51").Value = "  +0.37%  "

# The leading apostrophe above forces Excel to keep these numeric-looking
# strings (e.g. "211.44") as text instead of silently converting them to
# real numbers; ClearFormats() then drops the auto-applied Text number
# format so the cells keep their original default (unstyled) appearance.
$ws.Range("D5").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D29").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D47").ClearFormats()

